# Insert a new data row at row 37 (pushing the existing rows 37-54 down to
# 38-55), and populate it with the new record. All other rows keep their
# original values - they just move down one row, which EntireRow.Insert()
# does natively (shifts cell values + formatting along with the rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 37; existing data (and row 37's own formatting,
# used as the template) shifts down to row 38 onward.
$ws.Range("A37").EntireRow.Insert()

# Populate the newly inserted row 37 with the new record.
$ws.Range("A37").Value = 4
$ws.Range("B37").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C37").Value = "Los Lagos"
$ws.Range("D37").Value = 44510
$ws.Range("E37").Value = 10
$ws.Range("F37").Value = 100112026
$ws.Range("G37").Value = "Haba"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 120
$ws.Range("K37").Value = 10000
$ws.Range("L37").Value = 10000
$ws.Range("M37").Value = 10000
$ws.Range("N37").Value = "$/saco 25 kilos"
$ws.Range("O37").Value = "Región Metropolitana"
$ws.Range("P37").Value = 400
$ws.Range("Q37").Value = 25
$ws.Range("R37").Value = "Hortaliza"
